$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "col" column data (D1:D4), keeping D1's header style but clearing its value
$ws.Range("D1:D4").ClearContents()

# Row 3: replace "Gathersburg" data with "Charleston" data
$ws.Range("A3").Value = "Charleston"
$ws.Range("B3").Value = 38.3498
$ws.Range("C3").Value = -81.6326

# Row 4: replace "Charleston" data with "Morgantown" data
$ws.Range("A4").Value = "Morgantown"
$ws.Range("B4").Value = 39.6295
$ws.Range("C4").Value = -79.9559

# Row 5: add a new city "Bowie"
$ws.Range("A5").Value = "Bowie"
$ws.Range("B5").Value = 39.0068
$ws.Range("C5").Value = -76.7791

# G4: a cell with just whitespace text
$ws.Range("G4").Value = "   "

# Update the selection to match the recorded state
$ws.Range("J19").Select()
